$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRange, $text)
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.ClearFormats()
}

Set-TextValue $ws.Range("D2") '26.002.85'
Set-TextValue $ws.Range("E2") '  +0.53%  '

Set-TextValue $ws.Range("D3") '1.641.68'
Set-TextValue $ws.Range("E3") '  +0.64%  '

Set-TextValue $ws.Range("E4") '  +0.35%  '

Set-TextValue $ws.Range("D5") '215.96'
Set-TextValue $ws.Range("E5") '  +0.74%  '

Set-TextValue $ws.Range("E6") '  +0.31%  '

Set-TextValue $ws.Range("E8") '  +0.31%  '

Set-TextValue $ws.Range("E9") '  +1.05%  '

Set-TextValue $ws.Range("D10") '19.60'
Set-TextValue $ws.Range("E10") '  +0.31%  '

Set-TextValue $ws.Range("E11") '  +0.45%  '

Set-TextValue $ws.Range("D12") '1.869.43'
Set-TextValue $ws.Range("E12") '  +0.68%  '

Set-TextValue $ws.Range("D13") '1.670.08'
Set-TextValue $ws.Range("E13") '  +4.99%  '

Set-TextValue $ws.Range("E14") '  +0.29%  '

Set-TextValue $ws.Range("E15") '  -0.28%  '

Set-TextValue $ws.Range("D16") '0.0₃0763'
Set-TextValue $ws.Range("E16") '  +1.07%  '

Set-TextValue $ws.Range("D17") '63.33'
Set-TextValue $ws.Range("E17") '  +1.12%  '

Set-TextValue $ws.Range("D18") '26.101.57'
Set-TextValue $ws.Range("E18") '  +0.92%  '

Set-TextValue $ws.Range("E19") '  +0.37%  '

Set-TextValue $ws.Range("D20") '194.86'
Set-TextValue $ws.Range("E20") '  +0.86%  '

Set-TextValue $ws.Range("E21") '  -0.66%  '

Set-TextValue $ws.Range("D22") '9.91'
Set-TextValue $ws.Range("E22") '  +0.00%  '

Set-TextValue $ws.Range("E23") '  -0.78%  '

Set-TextValue $ws.Range("E24") '  +4.13%  '

Set-TextValue $ws.Range("E25") '  -1.65%  '

Set-TextValue $ws.Range("E26") '  +0.68%  '

Set-TextValue $ws.Range("D27") '143.05'
Set-TextValue $ws.Range("E27") '  -0.03%  '

Set-TextValue $ws.Range("D28") '6.88'
Set-TextValue $ws.Range("E28") '  +0.61%  '

Set-TextValue $ws.Range("E29") '  +0.54%  '

Set-TextValue $ws.Range("E30") '  +1.18%  '

Set-TextValue $ws.Range("D31") '0.0497'
Set-TextValue $ws.Range("E31") '  -0.07%  '

Set-TextValue $ws.Range("E32") '  +0.15%  '

Set-TextValue $ws.Range("E33") '  +1.22%  '

Set-TextValue $ws.Range("E34") '  -1.74%  '

Set-TextValue $ws.Range("E35") '  +1.71%  '

Set-TextValue $ws.Range("E36") '  +0.51%  '

Set-TextValue $ws.Range("D37") '1.127.78'
Set-TextValue $ws.Range("E37") '  -0.81%  '

Set-TextValue $ws.Range("D38") '0.541'
Set-TextValue $ws.Range("E38") '  -1.19%  '

Set-TextValue $ws.Range("E39") '  -0.23%  '

Set-TextValue $ws.Range("E40") '  +0.43%  '

Set-TextValue $ws.Range("E41") '  +0.43%  '

Set-TextValue $ws.Range("D42") '99.24'
Set-TextValue $ws.Range("E42") '  +0.18%  '

Set-TextValue $ws.Range("D44") '1.778.78'
Set-TextValue $ws.Range("E44") '  +0.72%  '

Set-TextValue $ws.Range("E45") '  +4.39%  '

Set-TextValue $ws.Range("D46") '56.63'
Set-TextValue $ws.Range("E46") '  +0.79%  '

Set-TextValue $ws.Range("D47") '0.0524'
Set-TextValue $ws.Range("E47") '  -0.61%  '

Set-TextValue $ws.Range("E48") '  +1.31%  '

Set-TextValue $ws.Range("E50") '  -0.24%  '

Set-TextValue $ws.Range("D51") '0.0955'
Set-TextValue $ws.Range("E51") '  -0.53%  '
